$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Résultats Algorithme")
$ws3 = $wb.Worksheets.Item("Resultats_merged")

# Update peak-detection / weight results on "Resultats_merged" (rows 2-31, cols K-S)
$ws3.Range("K2").Value = 290.60400000000004
$ws3.Range("L2").Value = 199
$ws3.Range("M2").Value = 26
$ws3.Range("N2").Value = 125.94099999999997
$ws3.Range("O2").Value = 5.476
$ws3.Range("P2").Value = 20.214999999999975
$ws3.Range("Q2").Value = 0.40399999999999636
$ws3.Range("R2").Value = 43.3
$ws3.Range("S2").Value = 23
$ws3.Range("K3").Value = 172.87400000000002
$ws3.Range("L3").Value = 191
$ws3.Range("M3").Value = 17
$ws3.Range("N3").Value = 63.284999999999968
$ws3.Range("O3").Value = 4.8680000000000003
$ws3.Range("P3").Value = 7.0779999999999745
$ws3.Range("Q3").Value = 0.20100000000002183
$ws3.Range("R3").Value = 36.6
$ws3.Range("S3").Value = 13
$ws3.Range("K5").Value = 272.40999999999997
$ws3.Range("L5").Value = 141
$ws3.Range("M5").Value = 44
$ws3.Range("N5").Value = 92.68900000000005
$ws3.Range("O5").Value = 2.5750000000000002
$ws3.Range("P5").Value = 8.8980000000000246
$ws3.Range("Q5").Value = 0.20199999999999818
$ws3.Range("R5").Value = 34
$ws3.Range("S5").Value = 36
$ws3.Range("K6").Value = 226.44399999999996
$ws3.Range("L6").Value = 169
$ws3.Range("M6").Value = 23
$ws3.Range("N6").Value = 90.380000000000223
$ws3.Range("O6").Value = 4.7569999999999997
$ws3.Range("P6").Value = 8.4920000000000755
$ws3.Range("Q6").Value = 0.20299999999997453
$ws3.Range("R6").Value = 39.9
$ws3.Range("S6").Value = 19
$ws3.Range("K7").Value = 490.779
$ws3.Range("L7").Value = 177
$ws3.Range("M7").Value = 34
$ws3.Range("N7").Value = 116.24700000000053
$ws3.Range("O7").Value = 5.0540000000000003
$ws3.Range("P7").Value = 14.154999999999973
$ws3.Range("Q7").Value = 0.20199999999999818
$ws3.Range("R7").Value = 23.7
$ws3.Range("S7").Value = 23
$ws3.Range("K8").Value = 288.45399999999995
$ws3.Range("L8").Value = 160
$ws3.Range("M8").Value = 39
$ws3.Range("N8").Value = 87.703999999999922
$ws3.Range("O8").Value = 2.7410000000000001
$ws3.Range("P8").Value = 8.4459999999999695
$ws3.Range("Q8").Value = 0.21600000000000819
$ws3.Range("R8").Value = 30.4
$ws3.Range("S8").Value = 32
$ws3.Range("K9").Value = 193.71000000000004
$ws3.Range("L9").Value = 178
$ws3.Range("M9").Value = 22
$ws3.Range("N9").Value = 51.976999999999748
$ws3.Range("O9").Value = 2.5990000000000002
$ws3.Range("P9").Value = 3.8979999999999109
$ws3.Range("Q9").Value = 0.21699999999998454
$ws3.Range("R9").Value = 26.8
$ws3.Range("S9").Value = 20
$ws3.Range("K10").Value = 51.322000000000003
$ws3.Range("L10").Value = 35
$ws3.Range("M10").Value = 10
$ws3.Range("N10").Value = 28.586000000000013
$ws3.Range("O10").Value = 3.573
$ws3.Range("P10").Value = 10.394999999999982
$ws3.Range("Q10").Value = 0.43299999999999272
$ws3.Range("R10").Value = 55.7
$ws3.Range("S10").Value = 8
$ws3.Range("K11").Value = 262.46799999999996
$ws3.Range("L11").Value = 146
$ws3.Range("M11").Value = 31
$ws3.Range("N11").Value = 106.54500000000004
$ws3.Range("O11").Value = 4.843
$ws3.Range("P11").Value = 8.4449999999999932
$ws3.Range("Q11").Value = 0.21600000000000819
$ws3.Range("R11").Value = 40.6
$ws3.Range("S11").Value = 22
$ws3.Range("K12").Value = 214.28099999999995
$ws3.Range("L12").Value = 113
$ws3.Range("M12").Value = 29
$ws3.Range("N12").Value = 81.538000000000238
$ws3.Range("O12").Value = 3.5449999999999999
$ws3.Range("P12").Value = 6.7119999999999891
$ws3.Range("Q12").Value = 0.64999999999997726
$ws3.Range("R12").Value = 38.1
$ws3.Range("S12").Value = 23
$ws3.Range("K13").Value = 395.53899999999999
$ws3.Range("L13").Value = 118
$ws3.Range("M13").Value = 45
$ws3.Range("N13").Value = 197.06600000000094
$ws3.Range("O13").Value = 7.2990000000000004
$ws3.Range("P13").Value = 14.076000000000022
$ws3.Range("Q13").Value = 0.43299999999999272
$ws3.Range("R13").Value = 49.8
$ws3.Range("S13").Value = 27
$ws3.Range("K14").Value = 250.34200000000001
$ws3.Range("L14").Value = 119
$ws3.Range("M14").Value = 24
$ws3.Range("N14").Value = 46.778999999999911
$ws3.Range("O14").Value = 2.2280000000000002
$ws3.Range("P14").Value = 3.6819999999999595
$ws3.Range("Q14").Value = 0.21600000000000819
$ws3.Range("R14").Value = 18.7
$ws3.Range("S14").Value = 21
$ws3.Range("K15").Value = 214.71500000000003
$ws3.Range("L15").Value = 131
$ws3.Range("M15").Value = 24
$ws3.Range("N15").Value = 90.717999999999734
$ws3.Range("O15").Value = 5.04
$ws3.Range("P15").Value = 9.3120000000000118
$ws3.Range("Q15").Value = 0.41699999999991633
$ws3.Range("R15").Value = 42.3
$ws3.Range("S15").Value = 18
$ws3.Range("K16").Value = 455.30900000000008
$ws3.Range("L16").Value = 178
$ws3.Range("M16").Value = 41
$ws3.Range("N16").Value = 204.19700000000046
$ws3.Range("O16").Value = 7.2930000000000001
$ws3.Range("P16").Value = 13.426000000000158
$ws3.Range("Q16").Value = 0.21600000000000819
$ws3.Range("R16").Value = 44.8
$ws3.Range("S16").Value = 28
$ws3.Range("K17").Value = 291.48700000000002
$ws3.Range("L17").Value = 106
$ws3.Range("M17").Value = 36
$ws3.Range("N17").Value = 91.391999999999882
$ws3.Range("O17").Value = 3.2639999999999998
$ws3.Range("P17").Value = 10.394999999999982
$ws3.Range("Q17").Value = 0.21699999999998454
$ws3.Range("R17").Value = 31.4
$ws3.Range("S17").Value = 28
$ws3.Range("K18").Value = 234.20500000000004
$ws3.Range("L18").Value = 126
$ws3.Range("M18").Value = 18
$ws3.Range("N18").Value = 84.456000000000017
$ws3.Range("O18").Value = 6.0330000000000004
$ws3.Range("P18").Value = 9.3129999999999882
$ws3.Range("Q18").Value = 0.21699999999998454
$ws3.Range("R18").Value = 36.1
$ws3.Range("S18").Value = 14
$ws3.Range("K19").Value = 436.46800000000007
$ws3.Range("L19").Value = 185
$ws3.Range("M19").Value = 58
$ws3.Range("N19").Value = 152.78499999999974
$ws3.Range("O19").Value = 3.056
$ws3.Range("P19").Value = 13.211000000000013
$ws3.Range("Q19").Value = 0.43299999999999272
$ws3.Range("R19").Value = 35
$ws3.Range("S19").Value = 50
$ws3.Range("K20").Value = 229.05300000000005
$ws3.Range("L20").Value = 131
$ws3.Range("M20").Value = 25
$ws3.Range("N20").Value = 127.07299999999987
$ws3.Range("O20").Value = 6.0510000000000002
$ws3.Range("P20").Value = 14.572000000000003
$ws3.Range("Q20").Value = 0.40399999999999636
$ws3.Range("R20").Value = 55.5
$ws3.Range("S20").Value = 21
$ws3.Range("K21").Value = 227.19399999999996
$ws3.Range("L21").Value = 139
$ws3.Range("M21").Value = 19
$ws3.Range("N21").Value = 123.40800000000013
$ws3.Range("O21").Value = 7.2590000000000003
$ws3.Range("P21").Value = 18.004999999999995
$ws3.Range("Q21").Value = 0.40399999999999636
$ws3.Range("R21").Value = 54.3
$ws3.Range("S21").Value = 17
$ws3.Range("K22").Value = 346.73500000000001
$ws3.Range("L22").Value = 207
$ws3.Range("M22").Value = 36
$ws3.Range("N22").Value = 180.65299999999979
$ws3.Range("O22").Value = 7.5270000000000001
$ws3.Range("P22").Value = 19.421000000000049
$ws3.Range("Q22").Value = 0.6069999999999709
$ws3.Range("R22").Value = 52.1
$ws3.Range("S22").Value = 24
$ws3.Range("K23").Value = 226.22100000000006
$ws3.Range("L23").Value = 132
$ws3.Range("M23").Value = 32
$ws3.Range("N23").Value = 116.14200000000011
$ws3.Range("O23").Value = 4.3019999999999996
$ws3.Range("P23").Value = 9.7100000000000364
$ws3.Range("Q23").Value = 0.6069999999999709
$ws3.Range("R23").Value = 51.3
$ws3.Range("S23").Value = 27
$ws3.Range("K24").Value = 240.54399999999998
$ws3.Range("L24").Value = 201
$ws3.Range("M24").Value = 31
$ws3.Range("N24").Value = 116.32499999999993
$ws3.Range("O24").Value = 4.1539999999999999
$ws3.Range("P24").Value = 13.149000000000001
$ws3.Range("Q24").Value = 0.20199999999999818
$ws3.Range("R24").Value = 48.4
$ws3.Range("S24").Value = 28
$ws3.Range("K25").Value = 198.25299999999993
$ws3.Range("L25").Value = 103
$ws3.Range("M25").Value = 23
$ws3.Range("N25").Value = 103.37600000000032
$ws3.Range("O25").Value = 5.1689999999999996
$ws3.Range("P25").Value = 14.36200000000008
$ws3.Range("Q25").Value = 0.20300000000020191
$ws3.Range("R25").Value = 52.1
$ws3.Range("S25").Value = 20
$ws3.Range("K26").Value = 199.298
$ws3.Range("L26").Value = 178
$ws3.Range("M26").Value = 23
$ws3.Range("N26").Value = 70.618000000000052
$ws3.Range("O26").Value = 3.5310000000000001
$ws3.Range("P26").Value = 9.964999999999975
$ws3.Range("Q26").Value = 0.64900000000000091
$ws3.Range("R26").Value = 35.4
$ws3.Range("S26").Value = 20
$ws3.Range("K27").Value = 213.16199999999992
$ws3.Range("L27").Value = 142
$ws3.Range("M27").Value = 20
$ws3.Range("N27").Value = 87.9549999999997
$ws3.Range("O27").Value = 5.8639999999999999
$ws3.Range("P27").Value = 10.182000000000016
$ws3.Range("Q27").Value = 0.2159999999998945
$ws3.Range("R27").Value = 41.3
$ws3.Range("S27").Value = 15
$ws3.Range("K28").Value = 260.60100000000011
$ws3.Range("L28").Value = 164
$ws3.Range("M28").Value = 22
$ws3.Range("N28").Value = 61.979000000000724
$ws3.Range("O28").Value = 4.7679999999999998
$ws3.Range("P28").Value = 7.1490000000001146
$ws3.Range("Q28").Value = 0.21600000000012187
$ws3.Range("R28").Value = 23.8
$ws3.Range("S28").Value = 13
$ws3.Range("K29").Value = 216.19400000000002
$ws3.Range("L29").Value = 168
$ws3.Range("M29").Value = 37
$ws3.Range("N29").Value = 91.061000000000263
$ws3.Range("O29").Value = 3.2519999999999998
$ws3.Range("P29").Value = 6.4990000000000236
$ws3.Range("Q29").Value = 0.21699999999998454
$ws3.Range("R29").Value = 42.1
$ws3.Range("S29").Value = 28
$ws3.Range("K30").Value = 207.74599999999998
$ws3.Range("L30").Value = 112
$ws3.Range("M30").Value = 24
$ws3.Range("N30").Value = 55.672000000000025
$ws3.Range("O30").Value = 3.2749999999999999
$ws3.Range("P30").Value = 5.6330000000000382
$ws3.Range("Q30").Value = 0.21600000000000819
$ws3.Range("R30").Value = 26.8
$ws3.Range("S30").Value = 17
$ws3.Range("K31").Value = 232.22400000000005
$ws3.Range("L31").Value = 85
$ws3.Range("M31").Value = 24
$ws3.Range("N31").Value = 59.657999999999447
$ws3.Range("O31").Value = 3.9769999999999999
$ws3.Range("P31").Value = 6.7149999999999181
$ws3.Range("Q31").Value = 0.21699999999998454
$ws3.Range("R31").Value = 25.7
$ws3.Range("S31").Value = 15

# Page setup for Resultats_merged sheet
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# Restore cursor / selection positions
$ws2.Activate()
$ws2.Range("B28").Select()
$ws3.Activate()
$ws3.Range("A8:XFD10").Select()
